$d = $word.ActiveDocument

# The document currently ends with an empty list paragraph (ilvl=0, numId=1).
# We fill it in with "Add bootstrap" and then append a run of new bullet
# points underneath it (sub-bullets at ilvl=1), followed by a new top level
# bullet "Material UI" with its own sub-bullet, finishing with a new empty
# trailing list paragraph (mirroring how the document ended before the edit).

$last = $d.Paragraphs.Last
$last.Range.Text = "Add bootstrap"

$last.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListIndent()
$p.Range.Text = "What’s bootstrap?"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Installing bootstrap."

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Ng bootstrap"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "Working with bootstrap components"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListOutdent()
$p.Range.Text = "Material UI"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.ListFormat.ListIndent()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
